$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # "Rubric"
$ws2 = $wb.Worksheets.Item(2)   # "Grade"

# ---------------------------------------------------------------
# Sheet "Rubric" (sheet1) cell content / value updates
# ---------------------------------------------------------------

# Row 10: model properties requirement text + points change, with text-wrap style
$ws1.Range("A10").WrapText = $true
$ws1.Range("A10").NumberFormat = "@"
$ws1.Range("A10").Value2 = ">= 5 properties including Date, AppUser and an int"
$ws1.Range("B10").Value2 = 3
$ws1.Rows.Item(10).RowHeight = 31.5

# Row 16: HTTP Post requirement text + points change, with text-wrap style
$ws1.Range("A16").WrapText = $true
$ws1.Range("A16").NumberFormat = "@"
$ws1.Range("A16").Value2 = "HTTP Post method to add date and echo entry"
$ws1.Range("B16").Value2 = 4
$ws1.Rows.Item(16).RowHeight = 31.5

# B3: "Possible" -> "Points" (done last so "Points" is appended after
# the two strings above in the shared-strings table, matching target order)
$ws1.Range("B3").Value2 = "Points"

# ---------------------------------------------------------------
# Sheet "Grade" (sheet2) cell content / value updates
# ---------------------------------------------------------------

# Row 11: model properties requirement text + points change, with text-wrap style
$ws2.Range("A11").WrapText = $true
$ws2.Range("A11").NumberFormat = "@"
$ws2.Range("A11").Value2 = ">= 5 properties including Date, AppUser and an int"
$ws2.Range("B11").Value2 = 3
$ws2.Range("C11").Value2 = 3
$ws2.Rows.Item(11).RowHeight = 31.5

# Row 17: HTTP Post requirement text + points change, with text-wrap style
$ws2.Range("A17").WrapText = $true
$ws2.Range("A17").NumberFormat = "@"
$ws2.Range("A17").Value2 = "HTTP Post method to add date and echo entry"
$ws2.Range("B17").Value2 = 4
$ws2.Range("C17").Value2 = 4
$ws2.Rows.Item(17).RowHeight = 31.5

# ---------------------------------------------------------------
# Column widths on "Rubric" sheet (best effort - COM ColumnWidth is
# quantized to character units, closest achievable values chosen)
# ---------------------------------------------------------------
$ws1.Columns.Item(1).ColumnWidth = 41.5
$ws1.Columns.Item(2).ColumnWidth = 5.166666666666666

# ---------------------------------------------------------------
# View state: work on the non-active sheet ("Grade") first, then
# activate "Rubric" last so it ends up as the selected/active tab,
# matching the original workbook's tab selection.
# ---------------------------------------------------------------
$ws2.Activate()
$win2 = $ws2.Parent.Windows.Item(1)
$win2.ScrollRow = 5
$win2.ScrollColumn = 1
$ws2.Range("A10:B21").Select()

$ws1.Activate()
$win1 = $ws1.Parent.Windows.Item(1)
$win1.Zoom = 120
$ws1.Range("C8").Select()
